# Update the bulk upload customer branch values (column A, rows 2-6) with
# parameterized Kolkata test-account values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Kol_10"
$ws.Range("A3").Value = "Kol_11"
$ws.Range("A4").Value = "Kol_12"
$ws.Range("A5").Value = "Kol_13"
$ws.Range("A6").Value = "Kol_14"

# Move the active selection to B10 (as captured in the saved view state).
$ws.Range("B10").Select()
